$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells so they match
# formatting (bold, border, centered) used by the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-46
$data = @(
    @(2, 2),
    @(5, 5),
    @(6, 6),
    @(6, 6),
    @(7, 8),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(4, 4),
    @(8, 8),
    @(4, 6),
    @(7, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(2, 6),
    @(7, 9),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 2),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 7),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 3),
    @(1, 6),
    @(1, 4),
    @(1, 2),
    @(1, 2)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
